$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 41.528285
$ws.Range("H2").Value = 124.584855
$ws.Range("I2").Value = 0.137866712381124
$ws.Range("J2").Value = 0.145939792231724
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8436446666666667
$ws.Range("N2").Value = 2.530934
$ws.Range("O2").Value = 0.038256548453167
$ws.Range("P2").Value = 0.03944386410459907
$ws.Range("Q2").Value = 35.03511615606334
$ws.Range("R2").Value = 315.31604540457
$ws.Range("S2").Value = 0.005274304562287307
$ws.Range("T2").Value = 0.005756429332241544

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 41.528285
$ws.Range("H3").Value = 124.584855
$ws.Range("I3").Value = 0.137866712381124
$ws.Range("J3").Value = 0.145939792231724
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.906580666666665
$ws.Range("N3").Value = 26.719742
$ws.Range("O3").Value = 0.4038845360958133
$ws.Range("P3").Value = 0.4164193425660044
$ws.Range("Q3").Value = 369.8750203008233
$ws.Range("R3").Value = 3328.87518270741
$ws.Range("S3").Value = 0.05568223317310517
$ws.Range("T3").Value = 0.06077215233535378

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 41.528285
$ws.Range("H4").Value = 124.584855
$ws.Range("I4").Value = 0.137866712381124
$ws.Range("J4").Value = 0.145939792231724
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.206402666666667
$ws.Range("N4").Value = 18.619208
$ws.Range("O4").Value = 0.2814402244434642
$ws.Range("P4").Value = 0.2901748959424717
$ws.Range("Q4").Value = 257.7412587660933
$ws.Range("R4").Value = 2319.67132889484
$ws.Range("S4").Value = 0.03880123847582605
$ws.Range("T4").Value = 0.04234806402470645

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 41.528285
$ws.Range("H5").Value = 124.584855
$ws.Range("I5").Value = 0.137866712381124
$ws.Range("J5").Value = 0.145939792231724
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.104251333333333
$ws.Range("N5").Value = 12.312754
$ws.Range("O5").Value = 0.1861144818446177
$ws.Range("P5").Value = 0.1918906599418865
$ws.Range("Q5").Value = 170.4425190822967
$ws.Range("R5").Value = 1533.98267174067
$ws.Range("S5").Value = 0.02565899173843383
$ws.Range("T5").Value = 0.02800448304312731

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 41.528285
$ws.Range("H6").Value = 124.584855
$ws.Range("I6").Value = 0.137866712381124
$ws.Range("J6").Value = 0.145939792231724
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.991415
$ws.Range("N6").Value = 3.98283
$ws.Range("O6").Value = 0.09030420916293774
$ws.Range("P6").Value = 0.06207123744503819
$ws.Range("Q6").Value = 82.700049673275
$ws.Range("R6").Value = 496.20029803965
$ws.Range("S6").Value = 0.0124499444314716
$ws.Range("T6").Value = 0.00905866349629488

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 72.02213166666667
$ws.Range("H7").Value = 216.066395
$ws.Range("I7").Value = 0.2391010009578718
$ws.Range("J7").Value = 0.2531020708300187
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8436446666666667
$ws.Range("N7").Value = 2.530934
$ws.Range("O7").Value = 0.038256548453167
$ws.Range("P7").Value = 0.03944386410459907
$ws.Range("Q7").Value = 60.76108726254778
$ws.Range("R7").Value = 546.8497853629301
$ws.Range("S7").Value = 0.009147179028345552
$ws.Range("T7").Value = 0.009983323686411864

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 72.02213166666667
$ws.Range("H8").Value = 216.066395
$ws.Range("I8").Value = 0.2391010009578718
$ws.Range("J8").Value = 0.2531020708300187
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.906580666666665
$ws.Range("N8").Value = 26.719742
$ws.Range("O8").Value = 0.4038845360958133
$ws.Range("P8").Value = 0.4164193425660044
$ws.Range("Q8").Value = 641.4709254744544
$ws.Range("R8").Value = 5773.23832927009
$ws.Range("S8").Value = 0.09656919685191465
$ws.Range("T8").Value = 0.1053965979371307

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 72.02213166666667
$ws.Range("H9").Value = 216.066395
$ws.Range("I9").Value = 0.2391010009578718
$ws.Range("J9").Value = 0.2531020708300187
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.206402666666667
$ws.Range("N9").Value = 18.619208
$ws.Range("O9").Value = 0.2814402244434642
$ws.Range("P9").Value = 0.2901748959424717
$ws.Range("Q9").Value = 446.9983500350178
$ws.Range("R9").Value = 4022.98515031516
$ws.Range("S9").Value = 0.06729263937424038
$ws.Range("T9").Value = 0.07344386706592478

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 72.02213166666667
$ws.Range("H10").Value = 216.066395
$ws.Range("I10").Value = 0.2391010009578718
$ws.Range("J10").Value = 0.2531020708300187
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.104251333333333
$ws.Range("N10").Value = 12.312754
$ws.Range("O10").Value = 0.1861144818446177
$ws.Range("P10").Value = 0.1918906599418865
$ws.Range("Q10").Value = 295.5969299224255
$ws.Range("R10").Value = 2660.37236930183
$ws.Range("S10").Value = 0.04450015890180375
$ws.Range("T10").Value = 0.04856792340423038

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 72.02213166666667
$ws.Range("H11").Value = 216.066395
$ws.Range("I11").Value = 0.2391010009578718
$ws.Range("J11").Value = 0.2531020708300187
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.991415
$ws.Range("N11").Value = 3.98283
$ws.Range("O11").Value = 0.09030420916293774
$ws.Range("P11").Value = 0.06207123744503819
$ws.Range("Q11").Value = 143.425953332975
$ws.Range("R11").Value = 860.55571999785
$ws.Range("S11").Value = 0.02159182680156743
$ws.Range("T11").Value = 0.01571035873632096

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 65.334877
$ws.Range("H12").Value = 196.004631
$ws.Range("I12").Value = 0.2169004738773853
$ws.Range("J12").Value = 0.2296015444621718
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.8436446666666667
$ws.Range("N12").Value = 2.530934
$ws.Range("O12").Value = 0.038256548453167
$ws.Range("P12").Value = 0.03944386410459907
$ws.Range("Q12").Value = 55.11942052837267
$ws.Range("R12").Value = 496.0747847553541
$ws.Range("S12").Value = 0.008297863488405072
$ws.Range("T12").Value = 0.009056372117971966

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 65.334877
$ws.Range("H13").Value = 196.004631
$ws.Range("I13").Value = 0.2169004738773853
$ws.Range("J13").Value = 0.2296015444621718
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 8.906580666666665
$ws.Range("N13").Value = 26.719742
$ws.Range("O13").Value = 0.4038845360958133
$ws.Range("P13").Value = 0.4164193425660044
$ws.Range("Q13").Value = 581.9103523472446
$ws.Range("R13").Value = 5237.193171125202
$ws.Range("S13").Value = 0.08760274727092981
$ws.Range("T13").Value = 0.09561052419707683

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 65.334877
$ws.Range("H14").Value = 196.004631
$ws.Range("I14").Value = 0.2169004738773853
$ws.Range("J14").Value = 0.2296015444621718
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 6.206402666666667
$ws.Range("N14").Value = 18.619208
$ws.Range("O14").Value = 0.2814402244434642
$ws.Range("P14").Value = 0.2901748959424717
$ws.Range("Q14").Value = 405.4945548391387
$ws.Range("R14").Value = 3649.450993552248
$ws.Range("S14").Value = 0.06104451804994505
$ws.Range("T14").Value = 0.0666246042725415

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 65.334877
$ws.Range("H15").Value = 196.004631
$ws.Range("I15").Value = 0.2169004738773853
$ws.Range("J15").Value = 0.2296015444621718
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 4.104251333333333
$ws.Range("N15").Value = 12.312754
$ws.Range("O15").Value = 0.1861144818446177
$ws.Range("P15").Value = 0.1918906599418865
$ws.Range("Q15").Value = 268.1507560404193
$ws.Range("R15").Value = 2413.356804363774
$ws.Range("S15").Value = 0.04036831930754161
$ws.Range("T15").Value = 0.04405839189052254

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 65.334877
$ws.Range("H16").Value = 196.004631
$ws.Range("I16").Value = 0.2169004738773853
$ws.Range("J16").Value = 0.2296015444621718
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.991415
$ws.Range("N16").Value = 3.98283
$ws.Range("O16").Value = 0.09030420916293774
$ws.Range("P16").Value = 0.06207123744503819
$ws.Range("Q16").Value = 130.108854080955
$ws.Range("R16").Value = 780.65312448573
$ws.Range("S16").Value = 0.01958702576056371
$ws.Range("T16").Value = 0.01425165198405896

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 72.34659833333335
$ws.Range("H17").Value = 217.039795
$ws.Range("I17").Value = 0.2401781740848285
$ws.Range("J17").Value = 0.2542423201304522
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.8436446666666667
$ws.Range("N17").Value = 2.530934
$ws.Range("O17").Value = 0.038256548453167
$ws.Range("P17").Value = 0.03944386410459907
$ws.Range("Q17").Value = 61.03482183539224
$ws.Range("R17").Value = 549.3133965185301
$ws.Range("S17").Value = 0.009188387954269418
$ws.Range("T17").Value = 0.01002829952486353

$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 72.34659833333335
$ws.Range("H18").Value = 217.039795
$ws.Range("I18").Value = 0.2401781740848285
$ws.Range("J18").Value = 0.2542423201304522
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 8.906580666666665
$ws.Range("N18").Value = 26.719742
$ws.Range("O18").Value = 0.4038845360958133
$ws.Range("P18").Value = 0.4164193425660044
$ws.Range("Q18").Value = 644.3608140147655
$ws.Range("R18").Value = 5799.24732613289
$ws.Range("S18").Value = 0.09700425042059042
$ws.Range("T18").Value = 0.1058714198011785

$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 72.34659833333335
$ws.Range("H19").Value = 217.039795
$ws.Range("I19").Value = 0.2401781740848285
$ws.Range("J19").Value = 0.2542423201304522
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 6.206402666666667
$ws.Range("N19").Value = 18.619208
$ws.Range("O19").Value = 0.2814402244434642
$ws.Range("P19").Value = 0.2901748959424717
$ws.Range("Q19").Value = 449.0121208202623
$ws.Range("R19").Value = 4041.10908738236
$ws.Range("S19").Value = 0.06759579922085554
$ws.Range("T19").Value = 0.07377473878802655

$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 72.34659833333335
$ws.Range("H20").Value = 217.039795
$ws.Range("I20").Value = 0.2401781740848285
$ws.Range("J20").Value = 0.2542423201304522
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 4.104251333333333
$ws.Range("N20").Value = 12.312754
$ws.Range("O20").Value = 0.1861144818446177
$ws.Range("P20").Value = 0.1918906599418865
$ws.Range("Q20").Value = 296.9286226717145
$ws.Range("R20").Value = 2672.35760404543
$ws.Range("S20").Value = 0.04470063642018424
$ws.Range("T20").Value = 0.04878672659498885

$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 72.34659833333335
$ws.Range("H21").Value = 217.039795
$ws.Range("I21").Value = 0.2401781740848285
$ws.Range("J21").Value = 0.2542423201304522
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 1.991415
$ws.Range("N21").Value = 3.98283
$ws.Range("O21").Value = 0.09030420916293774
$ws.Range("P21").Value = 0.06207123744503819
$ws.Range("Q21").Value = 144.072101119975
$ws.Range("R21").Value = 864.4326067198501
$ws.Range("S21").Value = 0.02168910006892882
$ws.Range("T21").Value = 0.01578113542139471

$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 49.988644
$ws.Range("H22").Value = 99.977288
$ws.Range("I22").Value = 0.1659536386987904
$ws.Range("J22").Value = 0.1171142723456333
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 0.8436446666666667
$ws.Range("N22").Value = 2.530934
$ws.Range("O22").Value = 0.038256548453167
$ws.Range("P22").Value = 0.03944386410459907
$ws.Range("Q22").Value = 42.17265290449867
$ws.Range("R22").Value = 253.035917426992
$ws.Range("S22").Value = 0.006348813419859646
$ws.Range("T22").Value = 0.004619439443110164

$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 49.988644
$ws.Range("H23").Value = 99.977288
$ws.Range("I23").Value = 0.1659536386987904
$ws.Range("J23").Value = 0.1171142723456333
$ws.Range("K23").Value = 3
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 8.906580666666665
$ws.Range("N23").Value = 26.719742
$ws.Range("O23").Value = 0.4038845360958133
$ws.Range("P23").Value = 0.4164193425660044
$ws.Range("Q23").Value = 445.2278902032826
$ws.Range("R23").Value = 2671.367341219696
$ws.Range("S23").Value = 0.06702610837927317
$ws.Range("T23").Value = 0.0487686482952646

$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 49.988644
$ws.Range("H24").Value = 99.977288
$ws.Range("I24").Value = 0.1659536386987904
$ws.Range("J24").Value = 0.1171142723456333
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 6.206402666666667
$ws.Range("N24").Value = 18.619208
$ws.Range("O24").Value = 0.2814402244434642
$ws.Range("P24").Value = 0.2901748959424717
$ws.Range("Q24").Value = 310.2496534246507
$ws.Range("R24").Value = 1861.497920547904
$ws.Range("S24").Value = 0.04670602932259714
$ws.Range("T24").Value = 0.03398362179127243

$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 49.988644
$ws.Range("H25").Value = 99.977288
$ws.Range("I25").Value = 0.1659536386987904
$ws.Range("J25").Value = 0.1171142723456333
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 4.104251333333333
$ws.Range("N25").Value = 12.312754
$ws.Range("O25").Value = 0.1861144818446177
$ws.Range("P25").Value = 0.1918906599418865
$ws.Range("Q25").Value = 205.1659587885253
$ws.Range("R25").Value = 1230.995752731152
$ws.Range("S25").Value = 0.03088637547665428
$ws.Range("T25").Value = 0.0224731350090174

$ws.Range("E26").Value = 2
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 49.988644
$ws.Range("H26").Value = 99.977288
$ws.Range("I26").Value = 0.1659536386987904
$ws.Range("J26").Value = 0.1171142723456333
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 1.991415
$ws.Range("N26").Value = 3.98283
$ws.Range("O26").Value = 0.09030420916293774
$ws.Range("P26").Value = 0.06207123744503819
$ws.Range("Q26").Value = 99.54813549126
$ws.Range("R26").Value = 398.19254196504
$ws.Range("S26").Value = 0.01498631210040617
$ws.Range("T26").Value = 0.007269427806968673

